$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F13").Value = "Intro+7+5"
$ws.Range("G13").Value = "15812__collingridge_night_near_windsor_chp13of13_seed15812.html"

[void]$ws.Range("G13").Select()
